# Reorder "Recorded By" (column G) values so that "System" (exact case)
# appears first in the comma-separated list, preserving the relative
# order of the remaining entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $value = $cell.Value2

    if ($value -ne $null -and $value -is [string] -and $value.Contains(",")) {
        $parts = $value -split ",\s*"
        if ($parts.Length -gt 1) {
            $idx = -1
            for ($i = 0; $i -lt $parts.Length; $i++) {
                if ($parts[$i].Equals("System")) {
                    $idx = $i
                    break
                }
            }
            if ($idx -gt 0) {
                $newValue = "System"
                for ($i = 0; $i -lt $parts.Length; $i++) {
                    if ($i -ne $idx) {
                        $newValue = $newValue + ", " + $parts[$i]
                    }
                }
                $cell.Value = $newValue
            }
        }
    }
}
